{"js": "// Add a new sentence (\"\u4f7f\u7528Git\u5206\u652f\u7b80\u5355\u53c8\u5feb\u6377\u3002\") as its own run at the end\n// of the \"\u591a\u4e91\u8f6c\u5c0f\u96e8\u2026\" diary paragraph \u2014 mirrors what Word itself does when\n// you click at the end of that paragraph and type the new sentence: the\n// typed text becomes a new <w:r>, and the paragraph-mark's pending run\n// formatting (<w:pPr><w:rPr>\u2026</w:rPr></w:pPr>) is consumed/dropped because\n// it is now redundant with the trailing run's own <w:rPr>.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst needle = \"\u591a\u4e91\u8f6c\u5c0f\u96e8\";\nconst newSentence = \"\u4f7f\u7528Git\u5206\u652f\u7b80\u5355\u53c8\u5feb\u6377\u3002\";\n\n// Locate the target paragraph by its distinctive leading text.\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nconst target = paragraphs.items.find((p) => p.text.indexOf(needle) !== -1);\nif (!target) {\n  throw new Error('Paragraph containing \"' + needle + '\" was not found.');\n}\n\n// Pull the paragraph's own canonical OOXML so we can reuse its exact\n// attributes (w14:paraId, rsids, \u2026) and its existing run(s) verbatim.\nconst ooxmlResult = target.getOoxml();\nawait context.sync();\n\nconst fullXml = ooxmlResult.value;\nconst pStart = fullXml.indexOf(\"<w:p \");\nconst sectStart = fullXml.indexOf(\"<w:sectPr\");\nconst pEndTag = \"</w:p>\";\nconst pEnd = fullXml.lastIndexOf(pEndTag, sectStart === -1 ? fullXml.length : sectStart) + pEndTag.length;\nif (pStart === -1 || pEnd === -1) {\n  throw new Error(\"Could not isolate the target paragraph's OOXML.\");\n}\nlet paragraphXml = fullXml.substring(pStart, pEnd);\n\n// Drop the paragraph-mark run-properties block: once a real trailing run\n// carries the formatting, Word no longer needs the pending mark override.\nparagraphXml = paragraphXml.replace(/<w:pPr>[\\s\\S]*?<\\/w:pPr>/, \"\");\n\n// Append the new sentence as its own run, matching the existing run's\n// formatting (an eastAsia font hint).\nconst newRunXml =\n  '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>' +\n  newSentence +\n  \"</w:t></w:r>\";\nparagraphXml = paragraphXml.replace(/<\\/w:p>$/, newRunXml + \"</w:p>\");\n\nconst packageXml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n  \"<w:body>\" +\n  paragraphXml +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\ntarget.getRange(\"Whole\").insertOoxml(packageXml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Add a new sentence (\"\u4f7f\u7528Git\u5206\u652f\u7b80\u5355\u53c8\u5feb\u6377\u3002\") as its own run at the end\n# of the \"\u591a\u4e91\u8f6c\u5c0f\u96e8\u2026\" diary paragraph \u2014 mirrors what Word itself does when\n# you click at the end of that paragraph and type the new sentence: the\n# typed text becomes a new <w:r>, and the paragraph-mark's pending run\n# formatting (<w:pPr><w:rPr>\u2026</w:rPr></w:pPr>) is consumed/dropped because\n# it is now redundant with the trailing run's own <w:rPr>.\n\n$d = $word.ActiveDocument\n\n$needle = \"\u591a\u4e91\u8f6c\u5c0f\u96e8\"\n$newSentence = \"\u4f7f\u7528Git\u5206\u652f\u7b80\u5355\u53c8\u5feb\u6377\u3002\"\n\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.StartsWith($needle)) {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Paragraph containing '$needle' was not found.\"\n}\n\n$rng = $target.Range\n\n# Pull the paragraph's own canonical OOXML so we can reuse its exact\n# attributes (w14:paraId, rsids, ...) and its existing run(s) verbatim.\n$xml = $rng.WordOpenXML\n\n$pStart = $xml.IndexOf(\"<w:p \")\n$closeTag = \"</w:p>\"\n$pEndFound = $xml.IndexOf($closeTag, $pStart)\nif ($pStart -lt 0 -or $pEndFound -lt 0) {\n    throw \"Could not isolate the target paragraph's OOXML.\"\n}\n$pEnd = $pEndFound + $closeTag.Length\n$paraXml = $xml.Substring($pStart, $pEnd - $pStart)\n\n# Drop the paragraph-mark run-properties block: once a real trailing run\n# carries the formatting, Word no longer needs the pending mark override.\n$paraXml = [regex]::Replace($paraXml, \"<w:pPr>.*?</w:pPr>\", \"\")\n\n# Append the new sentence as its own run, matching the existing run's\n# formatting (an eastAsia font hint).\n$newRunXml = \"<w:r><w:rPr><w:rFonts w:hint=`\"eastAsia`\"/></w:rPr><w:t>\" + $newSentence + \"</w:t></w:r>\"\n$paraXml = $paraXml -replace \"</w:p>$\", ($newRunXml + \"</w:p>\")\n\n$pkgXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n    '<w:body>' + $paraXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n\n$rng.InsertXML($pkgXml)\n"}
